$wb = $excel.ActiveWorkbook

# --- Sheet "os": F3:F12 = "auto", plus J5 and J8 = "auto" ---
$wsOs = $wb.Worksheets.Item("os")
$wsOs.Range("F3").Value = "auto"
$wsOs.Range("F4").Value = "auto"
$wsOs.Range("F5").Value = "auto"
$wsOs.Range("J5").Value = "auto"
$wsOs.Range("F6").Value = "auto"
$wsOs.Range("F7").Value = "auto"
$wsOs.Range("F8").Value = "auto"
$wsOs.Range("J8").Value = "auto"
$wsOs.Range("F9").Value = "auto"
$wsOs.Range("F10").Value = "auto"
$wsOs.Range("F11").Value = "auto"
$wsOs.Range("F12").Value = "auto"
$wsOs.Activate()
$wsOs.Range("F4:F12").Select()

# --- Sheet "time": F3:F11 = "auto" ---
$wsTime = $wb.Worksheets.Item("time")
$wsTime.Range("F3").Value = "auto"
$wsTime.Range("F4").Value = "auto"
$wsTime.Range("F5").Value = "auto"
$wsTime.Range("F6").Value = "auto"
$wsTime.Range("F7").Value = "auto"
$wsTime.Range("F8").Value = "auto"
$wsTime.Range("F9").Value = "auto"
$wsTime.Range("F10").Value = "auto"
$wsTime.Range("F11").Value = "auto"
$wsTime.Activate()
$wsTime.Range("F3:F11").Select()

# --- Sheet "device": F3:F10 = "auto", plus J6 = "auto" and J7 = 2 ---
$wsDevice = $wb.Worksheets.Item("device")
$wsDevice.Range("F3").Value = "auto"
$wsDevice.Range("F4").Value = "auto"
$wsDevice.Range("F5").Value = "auto"
$wsDevice.Range("F6").Value = "auto"
$wsDevice.Range("J6").Value = "auto"
$wsDevice.Range("F7").Value = "auto"
$wsDevice.Range("J7").Value = 2
$wsDevice.Range("F8").Value = "auto"
$wsDevice.Range("F9").Value = "auto"
$wsDevice.Range("F10").Value = "auto"
$wsDevice.Activate()
$wsDevice.Range("F3:F10").Select()

# --- Sheet "modbus": F3 = "auto", J3 = 4 ---
$wsModbus = $wb.Worksheets.Item("modbus")
$wsModbus.Range("F3").Value = "auto"
$wsModbus.Range("J3").Value = 4
$wsModbus.Activate()
$wsModbus.Range("F9").Select()
